$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the bibliography entry for Subramonyam et al. (2024). There are two
# paragraphs that start with this citation text in the document (one in the
# "First Paragraph" styled intro, one in the "Bibliography" list) -- we only
# want the Bibliography list entry that still has the old
# "(arXiv:2309.14459). arXiv." tail.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Style.NameLocal -eq "Bibliography") {
        $t = $para.Range.Text
        if ($t -like "*Subramonyam*" -and $t -like "*LLM Interfaces (arXiv:2309.14459). arXiv.*") {
            $target = $para
            break
        }
    }
}

if ($target -eq $null) {
    throw "Could not find the Subramonyam et al. Bibliography paragraph"
}

$enDash = [char]8211

# ---------------------------------------------------------------------------
# Step 1: swap the old "(arXiv:2309.14459). arXiv." tail (which sits right
# after the italic "LLM Interfaces") for the new journal citation text. This
# also removes the old italic run boundary so the whole run momentarily picks
# up the italic formatting of "LLM Interfaces" -- fixed up in steps 2-4.
# ---------------------------------------------------------------------------
$r1 = $target.Range
$replacement = "LLM Interfaces. Proceedings of the CHI Conference on Human Factors in Computing Systems, 1" + $enDash + "19."
$found1 = $r1.Find.Execute("LLM Interfaces (arXiv:2309.14459). arXiv.", $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)
if (-not $found1) {
    throw "Could not find the old citation tail to replace"
}

# ---------------------------------------------------------------------------
# Step 2: the title "Bridging the Gulf of Envisioning: Cognitive Design
# Challenges in LLM Interfaces." is no longer italicized (only the new
# journal name is).
# ---------------------------------------------------------------------------
$r2 = $target.Range
$found2 = $r2.Find.Execute("Bridging the Gulf of Envisioning: Cognitive Design Challenges in LLM Interfaces.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $r2.Italic = $false
} else {
    throw "Could not find the (now non-italic) title text"
}

# ---------------------------------------------------------------------------
# Step 3: italicize the new journal/proceedings name.
# ---------------------------------------------------------------------------
$r3 = $target.Range
$found3 = $r3.Find.Execute("Proceedings of the CHI Conference on Human Factors in Computing Systems", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $r3.Italic = $true
} else {
    throw "Could not find the new journal name text to italicize"
}

# ---------------------------------------------------------------------------
# Step 4: the trailing page range ", 1-19." is not italic.
# ---------------------------------------------------------------------------
$r4 = $target.Range
$found4 = $r4.Find.Execute(", 1" + $enDash + "19.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $r4.Italic = $false
} else {
    throw "Could not find the trailing page range text"
}

Write-Host "Updated:" $target.Range.Text
